# Applies: add a new worksheet "validLoginData" after "inValidLoginData",
# populate it with a small userName/password table, make it the active
# (selected) sheet/tab, and clear the previous selection/active-tab state
# on the first sheet (replacing its selection with A1:B1).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws2.Name = "validLoginData"

# --- Fill in the valid-login data table ---
$ws2.Range("A1").Value = "userName"
$ws2.Range("B1").Value = "password"
$ws2.Range("A2").Value = "Admin"
$ws2.Range("B2").Value = "admin123"

# --- Column widths (approx. 13.29 / 13.14 displayed Excel character widths) ---
$ws2.Columns.Item(1).ColumnWidth = 12.45
$ws2.Columns.Item(2).ColumnWidth = 12.3

# --- Header row formatting: bold, centered, wrapped, taller row ---
$ws2.Rows.Item(1).RowHeight = 30
$headerRange = $ws2.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$headerRange.WrapText = $true

# --- Reset the selection on the original sheet (no more tab-selected/activeCell) ---
$ws1.Activate()
$ws1.Range("A1:B1").Select()

# --- Make the new sheet the active / selected tab ---
$ws2.Activate()
$ws2.Range("B2").Select()
